$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 1057, 1057, 1057, 0.01029562950134277),
    @(1, 1040, 1040, 1040, 0.01057639122009277),
    @(2, 973,  973,  973,  0.01228516896565755),
    @(3, 1224, 1224, 1224, 0.01205418109893799),
    @(4, 883,  883,  883,  0.01211105982462565),
    @(5, 1040, 1040, 1040, 0.01241656939188639),
    @(6, 1053, 1053, 1053, 0.01220947901407878),
    @(7, 957,  957,  957,  0.01029446919759115),
    @(8, 886,  886,  886,  0.01220994790395101),
    @(9, 1049, 1049, 1049, 0.0117668628692627)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
